$wb = $excel.ActiveWorkbook

# Rename sheets with updated timestamp-based names
$wb.Worksheets.Item(1).Name = "GNG_TO-16502911885366724"
$wb.Worksheets.Item(2).Name = "NB_TO-16502911914245977"
$wb.Worksheets.Item(3).Name = "RS_TO-16502911914266021"
$wb.Worksheets.Item(4).Name = "TOL_TO-16502911914885993"
$wb.Worksheets.Item(5).Name = "vSAT_TO-1650291191583606"

# Sheet 1: GNG
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("B2").Value = "go_stims-16502911884996686.csv"
$ws1.Range("B3").Value = "GNG_stims-16502911885196724.csv"
$ws1.Range("B4").Value = "go_stims-16502911885216691.csv"
$ws1.Range("B5").Value = "GNG_stims-16502911885356708.csv"

# Sheet 2: NB
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("B2").Value = "OB-1650291189652669.csv"
$ws2.Range("B3").Value = "TB-16502911909765978.csv"
$ws2.Range("B4").Value = "ZB-match_0-16502911886716685.csv"
$ws2.Range("B5").Value = "OB-16502911893096697.csv"
$ws2.Range("B6").Value = "ZB-match_2-16502911887646728.csv"
$ws2.Range("B7").Value = "ZB-match_6-16502911891436703.csv"
$ws2.Range("B8").Value = "OB-1650291190665151.csv"
$ws2.Range("B9").Value = "TB-1650291191086608.csv"
$ws2.Range("B10").Value = "TB-16502911914116008.csv"

# Sheet 3: RS
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("B2").Value = "eyes open"
$ws3.Range("B3").Value = "eyes closed"

# Sheet 4: TOL
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("B2").Value = "MM_stims-1650291191440598.csv"
$ws4.Range("B3").Value = "ZM_stims-16502911914296005.csv"
$ws4.Range("B4").Value = "MM_stims-16502911914715955.csv"
$ws4.Range("B5").Value = "ZM_stims-16502911914425974.csv"
$ws4.Range("B6").Value = "MM_stims-16502911914876156.csv"
$ws4.Range("B7").Value = "ZM_stims-16502911914736.csv"

# Sheet 5: vSAT
$ws5 = $wb.Worksheets.Item(5)
$ws5.Range("B2").Value = "vSAT_stims-1650291191537609.csv"
$ws5.Range("B3").Value = "SAT_stims-1650291191494636.csv"
$ws5.Range("B4").Value = "vSAT_stims-16502911915675972.csv"
$ws5.Range("B5").Value = "SAT_stims-16502911915216074.csv"
